# Auto-generated edit script: updates cached market-price/profit values
# in each Leve-tracking worksheet, matching the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 333333340
$ws.Range("I8").Value = 333333340
$ws.Range("K8").Value = 1000000020
$ws.Range("M8").Value = -999999881
$ws.Range("H19").Value = 1258.5294
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 1258.5294
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 1258.5294
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -1608.5294
$ws.Range("H40").Value = 4674.3335
$ws.Range("I40").Value = 4572.4
$ws.Range("J40").Value = 4878.2
$ws.Range("K40").Value = 4572.4
$ws.Range("L40").Value = 4878.2
$ws.Range("M40").Value = -4397.4
$ws.Range("N40").Value = -5228.2
$ws.Range("H100").Value = 3977822.8
$ws.Range("J100").Value = 10156587
$ws.Range("L100").Value = 10156587
$ws.Range("N100").Value = -10157669
$ws.Range("H113").Value = 11461
$ws.Range("I113").Value = 15735.272
$ws.Range("J113").Value = 3624.8333
$ws.Range("K113").Value = 15735.272
$ws.Range("L113").Value = 3624.8333
$ws.Range("M113").Value = -12481.272
$ws.Range("N113").Value = -10132.8333
$ws.Range("H138").Value = 4587.9785
$ws.Range("J138").Value = 4911.7856
$ws.Range("L138").Value = 14735.3568
$ws.Range("N138").Value = -25015.3568

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 157675.36
$ws.Range("I45").Value = 270381.88
$ws.Range("K45").Value = 270381.88
$ws.Range("M45").Value = -270004.88
$ws.Range("H61").Value = 6313.154
$ws.Range("I61").Value = 6915.636
$ws.Range("K61").Value = 6915.636
$ws.Range("M61").Value = -6703.636
$ws.Range("H63").Value = 3022
$ws.Range("I63").Value = 1777.5
$ws.Range("J63").Value = 8000
$ws.Range("K63").Value = 1777.5
$ws.Range("L63").Value = 8000
$ws.Range("M63").Value = -1091.5
$ws.Range("N63").Value = -9372
$ws.Range("H66").Value = 3022
$ws.Range("I66").Value = 1777.5
$ws.Range("J66").Value = 8000
$ws.Range("K66").Value = 8887.5
$ws.Range("L66").Value = 40000
$ws.Range("M66").Value = -5455.5
$ws.Range("N66").Value = -46864
$ws.Range("H102").Value = 7053.0835
$ws.Range("I102").Value = 2874.625
$ws.Range("K102").Value = 2874.625
$ws.Range("M102").Value = -1252.625
$ws.Range("H132").Value = 2959.6758
$ws.Range("I132").Value = 1954.3334
$ws.Range("K132").Value = 5863.0002
$ws.Range("M132").Value = -3333.0002
$ws.Range("H136").Value = 6313.154
$ws.Range("I136").Value = 6915.636
$ws.Range("K136").Value = 20746.908
$ws.Range("M136").Value = -18196.908

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4255.4
$ws.Range("I107").Value = 4544
$ws.Range("J107").Value = 3461.75
$ws.Range("K107").Value = 4544
$ws.Range("L107").Value = 3461.75
$ws.Range("M107").Value = -2624
$ws.Range("N107").Value = -7301.75
$ws.Range("H134").Value = 2236.5312
$ws.Range("I134").Value = 1775.8572
$ws.Range("K134").Value = 5327.571599999999
$ws.Range("M134").Value = -2792.571599999999
$ws.Range("H140").Value = 89799
$ws.Range("J140").Value = 89799
$ws.Range("L140").Value = 89799
$ws.Range("N140").Value = -100159

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 11799.833
$ws.Range("I41").Value = 2749.75
$ws.Range("J41").Value = 29900
$ws.Range("K41").Value = 2749.75
$ws.Range("L41").Value = 29900
$ws.Range("M41").Value = -2321.75
$ws.Range("N41").Value = -30756
$ws.Range("H62").Value = 19500
$ws.Range("I62").Value = 19000
$ws.Range("K62").Value = 19000
$ws.Range("M62").Value = -18376
$ws.Range("H65").Value = 19500
$ws.Range("I65").Value = 19000
$ws.Range("K65").Value = 95000
$ws.Range("M65").Value = -91880
$ws.Range("H99").Value = 213456.7
$ws.Range("I99").Value = 421577
$ws.Range("K99").Value = 421577
$ws.Range("M99").Value = -420079
$ws.Range("H105").Value = 82480.46000000001
$ws.Range("I105").Value = 118156.5
$ws.Range("K105").Value = 118156.5
$ws.Range("M105").Value = -116409.5
$ws.Range("H122").Value = 725
$ws.Range("J122").Value = 1000
$ws.Range("L122").Value = 3000
$ws.Range("N122").Value = -7900
$ws.Range("H126").Value = 213456.7
$ws.Range("I126").Value = 421577
$ws.Range("K126").Value = 1264731
$ws.Range("M126").Value = -1262261
$ws.Range("H134").Value = 5220559
$ws.Range("I134").Value = 6959023.5
$ws.Range("J134").Value = 5166.3335
$ws.Range("K134").Value = 20877070.5
$ws.Range("L134").Value = 15499.0005
$ws.Range("M134").Value = -20874535.5
$ws.Range("N134").Value = -20569.0005

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 471
$ws.Range("I92").Value = 479.66666
$ws.Range("J92").Value = 458
$ws.Range("K92").Value = 1438.99998
$ws.Range("L92").Value = 1374
$ws.Range("M92").Value = -190.9999800000001
$ws.Range("N92").Value = -3870
$ws.Range("H113").Value = 2710.5715
$ws.Range("I113").Value = 849.5
$ws.Range("J113").Value = 3455
$ws.Range("K113").Value = 2548.5
$ws.Range("L113").Value = 10365
$ws.Range("M113").Value = -378.5
$ws.Range("N113").Value = -14705
$ws.Range("H132").Value = 46631.273
$ws.Range("I132").Value = 809.8
$ws.Range("K132").Value = 7288.2
$ws.Range("M132").Value = -4758.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 11435.553
$ws.Range("I122").Value = 10012.29
$ws.Range("K122").Value = 30036.87
$ws.Range("M122").Value = -27586.87
$ws.Range("H132").Value = 4115.067
$ws.Range("I132").Value = 2351.7
$ws.Range("K132").Value = 7055.099999999999
$ws.Range("M132").Value = -4525.099999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 42897.727
$ws.Range("J7").Value = 4833.3335
$ws.Range("L7").Value = 4833.3335
$ws.Range("N7").Value = -5057.3335
$ws.Range("H16").Value = 733.3333
$ws.Range("J16").Value = 400
$ws.Range("L16").Value = 400
$ws.Range("N16").Value = -740
$ws.Range("H46").Value = 3431.1765
$ws.Range("I46").Value = 1419.5714
$ws.Range("J46").Value = 4839.3
$ws.Range("K46").Value = 1419.5714
$ws.Range("L46").Value = 4839.3
$ws.Range("M46").Value = -1231.5714
$ws.Range("N46").Value = -5215.3
$ws.Range("H68").Value = 4632.143
$ws.Range("I68").Value = 3783.4443
$ws.Range("K68").Value = 3783.4443
$ws.Range("M68").Value = -3034.4443
$ws.Range("H71").Value = 4632.143
$ws.Range("I71").Value = 3783.4443
$ws.Range("K71").Value = 18917.2215
$ws.Range("M71").Value = -15173.2215
$ws.Range("H126").Value = 42897.727
$ws.Range("J126").Value = 4833.3335
$ws.Range("L126").Value = 14500.0005
$ws.Range("N126").Value = -19440.0005
$ws.Range("H136").Value = 10474.913
$ws.Range("I136").Value = 10071.77
$ws.Range("K136").Value = 30215.31
$ws.Range("M136").Value = -27665.31

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 224321.08
$ws.Range("J62").Value = 3681.182
$ws.Range("L62").Value = 3681.182
$ws.Range("N62").Value = -4929.182
$ws.Range("H65").Value = 224321.08
$ws.Range("J65").Value = 3681.182
$ws.Range("L65").Value = 18405.91
$ws.Range("N65").Value = -24645.91
$ws.Range("H132").Value = 15808.333
$ws.Range("I132").Value = 17632.217
$ws.Range("K132").Value = 52896.651
$ws.Range("M132").Value = -50366.651
